# Update test data in the "IAM" worksheet: the request-body JSON strings in
# column H for several rows gained an extra "app":"cmty" field (one of them
# malformed, matching the source data exactly).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAM")

$ws.Range("H13").Value = '{"loginid":"(OPQA-542_email)","password":"1Platform!","app":"cmty"}'
$ws.Range("H14").Value = '{"loginid":"(OPQA-542_email)","password":"1Platform","app":"cmty"}'
$ws.Range("H21").Value = '{"loginid":"(OPQA-542_email)","password":"1Platform!","app":"cmty"}'
$ws.Range("H23").Value = '{"loginid":"(OPQA-542_email)","password":"1Platform!""app":"cmty"}'
$ws.Range("H27").Value = '{"loginid":"(OPQA-542_email)","password":"1Platform!","app":"cmty"}'
$ws.Range("H30").Value = '{"loginid":"(ddMMMyyyy_HHmmss)@sharklasers.com","password":"Neon@123",app="cmty"}'
$ws.Range("H33").Value = '{"loginid":"neontestuser007+(ddMMMyyyy_HHmmss)@gmail.com","password":"Neon@123","app":"cmty"}'
$ws.Range("H62").Value = '{"loginid":"project.neon2@gmail.com","password":"1Platform!","app":"cmty"}'

# Restore the sheet's active selection to match where the author left off
# (scrolled/selected near the bottom of the data, around A61/A62).
$ws.Activate()
$ws.Range("A62").Select()
